# Applies the "Entrega final - Reto 1" edit:
#  1. Collapses every "O(NlogN)" run-split (with spell-check proofErr markers
#     around "NlogN") into a single clean run, in all 4 occurrences.
#  2. Re-saves the three long "Encontramos la complejidad..." analysis
#     paragraphs (Requerimiento 1, 2 and 4) with identical text so Word
#     coalesces the runs that used to be split up for spell-check purposes.
#  3. Inserts the brand-new "Primero encontramos..." analysis paragraph
#     (plus an extra blank paragraph) for Requerimiento 3, right after its
#     "O(NlogN)" line.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) O(NlogN) clean-up, everywhere in the document.
# ---------------------------------------------------------------------
$target = "O(NlogN)"
$d.Content.Find.Execute($target, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $target, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Re-assert the unchanged long paragraphs so the engine coalesces runs.
# ---------------------------------------------------------------------
$req1Text = "Encontramos la complejidad de N log N en el algoritmo de ordenamiento de tipo Merge que se utiliza al inicio del código para ordenar por países. Este ordenamiento está realizado con el objetivo de realizar una futura búsqueda binaria (complejidad de log(N)) para encontrar el país que ingresa el usuario por parámetro. Luego de la búsqueda binaria se realiza una sublista que contiene solamente los videos pertenecientes al país de interés. Una vez creada esta sublista se procede a ordenar esta nueva lista por orden de categorías. Este ordenamiento también tiene un orden de complejidad (n log n) no obstante, toca considerar que este n es menor al N inicial ya que estamos trabajando con una sublista de la inicial. Una vez ordenada esta última lista de países se procede a realizar nuevamente una búsqueda binaria para encontrar la categoría que se ingresó por parámetro para así poder tener acceso nuevamente a los videos que pertenecen a la categoría de interés para una vez más realizar un nuevo ordenamiento por vistas para dar el resultado del primer requerimiento."

$req2Text = "Encontramos la complejidad de N log N en el algoritmo de ordenamiento de tipo Merge que se utiliza al inicio del código para ordenar por países. Este ordenamiento está realizado con el objetivo de realizar una futura búsqueda binaria (complejidad de log(N)) para encontrar el país que ingresa el usuario por parámetro. Luego de la búsqueda binaria se realiza una sublista que contiene solamente los videos pertenecientes al país de interés. Una vez creada esta sublista se procede a ordenar esta nueva lista por id. Este ordenamiento también tiene un orden de complejidad (n log n) no obstante, toca considerar que este n es menor al N inicial ya que estamos trabajando con una sublista de la inicial. Una vez ordenada esta última lista de países se procede a realizar un contento de orden de complejidad n donde n es la longitud total de la lista ordenada por id del país de interés para poder encontrar cual es el video que más veces ha sido tendencia en ese país. Una vez se concluye el conteo y se devuelve cual es video que más ha aparecido en la lista se realiza una búsqueda binaria para poder acceder a todos los datos del video para poder imprimir todos los datos que se necesitan la interfaz."

$req4Text = "Encontramos la complejidad de N log N en el algoritmo de ordenamiento de tipo Merge que se utiliza al inicio del código para ordenar por países. Este ordenamiento está realizado con el objetivo de realizar una futura búsqueda binaria (complejidad de log(N)) para encontrar el país que ingresa el usuario por parámetro. Luego de la búsqueda binaria se realiza una sublista que contiene solamente los videos pertenecientes al país de interés. Una vez creada esta sublista se procede a ordenar esta nueva lista por likes. Este ordenamiento también tiene un orden de complejidad (n log n) no obstante, toca considerar que este n es menor al N inicial ya que estamos trabajando con una sublista de la inicial. Una vez terminado este ordenamiento creamos una nueva lista de tipo arraylist debido a su utilidad para nuestro caso de estudio. En esta nueva lista se rellenará con los videos que cumplan con los tags de interés y además se verifica que los videos ya no existan en la lista para evitar repeticiones. Una vez realizado esto se devuelven los resultados y se imprimen en view.py"

foreach ($txt in @($req1Text, $req2Text, $req4Text)) {
    $rng = $d.Content
    $rng.Find.Execute($txt, $false, $false, $false, $false, $false, `
                       $true, 1, $false, $txt, 2) | Out-Null
}

# ---------------------------------------------------------------------
# 3) Insert the new Requerimiento 3 analysis paragraph.
# ---------------------------------------------------------------------
# Locate the (now single-run) "O(NlogN)" paragraph that belongs to
# Requerimiento 3: it is the one directly followed by an empty paragraph
# and then the "Requerimiento 4:" paragraph.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "O(NlogN)") {
        $nextPara = $p.Next()
        if ($nextPara -ne $null -and $nextPara.Range.Text.TrimEnd("`r") -eq "") {
            $afterNext = $nextPara.Next()
            if ($afterNext -ne $null -and $afterNext.Range.Text.TrimEnd("`r") -eq "Requerimiento 4:") {
                $targetPara = $p
            }
        }
    }
}

$newText = "Primero encontramos la complejidad de N log N en el algoritmo de ordenamiento de tipo Merge que se utiliza al inicio del código para ordenar una copia del catálogo por categorías. Este ordenamiento está realizado con el objetivo de hacer una búsqueda binaria de complejidad de log(N) para encontrar la categoría que fue indicada por el usuario. Luego de la búsqueda binaria se crea una sublista que contiene solamente los videos pertenecientes a la categoría indicada. Una vez creada esta sublista se procede a ordenarla teniendo en cuenta los video id’s, ordenamiento que también tiene un orden de complejidad (n log n) aunque en este caso se están manejando menos datos ya que se está trabajando con una sublista de la inicial. Luego se realiza otro ordenamiento de tipo merge (complejidad n log n) partiendo de la sublista mencionada previamente para así obtener una lista ordenada por ambos factores (id y fecha de trending). Una vez ordenada se procede a realizar un contento de orden de complejidad n (longitud de la lista ordenada) para poder encontrar cual es el video que más veces ha sido tendencia considerando que aunque sea trending en varios países el mismo día estos se deben contar una sola vez. Cuando se concluye el conteo y se devuelve cual es el video que más ha aparecido en la lista se realiza una búsqueda binaria para poder acceder a los datos del video para poder imprimir en consola los datos solicitados."

if ($targetPara -ne $null) {
    $insertRange = $targetPara.Range
    $insertRange.InsertParagraphAfter()

    $newPara = $targetPara.Next()
    $newPara.Range.InsertAfter($newText)

    # Extra blank paragraph, matching the source edit.
    $newPara.Range.InsertParagraphAfter()
}
